$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.768.98'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.634.75'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.79'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.499'
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.56'
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.860.80'
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.637.63'
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.562'
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.19'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.810.87'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.45'
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.38'
$ws.Range("E21").Value = '  -0.79%  '
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.30'
$ws.Range("E23").Value = '  +2.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.83'
$ws.Range("E24").Value = '  +3.50%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.46'
$ws.Range("E26").Value = '  +1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.123'
$ws.Range("E27").Value = '  +1.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.89'
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.48'
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("E33").Value = '  -0.74%  '
$ws.Range("E34").Value = '  -0.85%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.129.09'
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.51'
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.56'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.65'
$ws.Range("E43").Value = '  +1.20%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.801'
$ws.Range("E44").Value = '  +0.48%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.769.99'
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.33'
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.417'
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0502'
$ws.Range("E48").Value = '  -0.23%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.43'
$ws.Range("E49").Value = '  +4.08%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.47'
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  -0.15%  '
